# Commit volunteer information and sign list
# Adds a new volunteer (Sanah Khan) row to the sign-up sheet and
# updates the workbook's background theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append new volunteer row (row 36) ---
$ws.Range("A36").Value = "Sanah Khan"
$ws.Range("B36").Value = "sanah.khan@mail.utoronto.ca"
$ws.Range("D36").Value = "UTM"
$ws.Range("E36").Value = "905-821-9984"

# Keep the active selection on the newly added row, same as the source file.
$ws.Range("E36").Select() | Out-Null

# --- Update theme "Background 1 / Light 1" color ---
# RGB(0x7D, 0xE1, 0x7D) == 7DE17D
$wb.Theme.ThemeColorScheme.Item(2).RGB = 8249725
